# Apply name-list corrections to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "黃錦萍"
$ws.Range("D2").Value = "葉雪容"
$ws.Range("E2").Value = "ANI"
$ws.Range("D3").Value = "葉紅志"
$ws.Range("D4").Value = "廖練雲"
$ws.Range("D5").Value = "廖苑雲"
$ws.Range("D6").Value = "廖江真"
$ws.Range("D7").Value = "廖睦堯"
$ws.Range("C8").Value = "張志謀"
$ws.Range("D8").Value = "廖維華"
$ws.Range("C9").Value = "張志謀"
$ws.Range("C10").Value = "廖鋼基"
$ws.Range("D10").Value = "廖富盛"
$ws.Range("C11").Value = "廖偉良"
$ws.Range("C12").Value = "廖素琼"
$ws.Range("D13").ClearContents()

$null = $ws.Range("D14").Select()
